$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 2411.6072
$ws.Range("I19").Value = 7382.857
$ws.Range("J19").Value = 754.5238000000001
$ws.Range("K19").Value = 7382.857
$ws.Range("L19").Value = 754.5238000000001
$ws.Range("M19").Value = -7207.857
$ws.Range("N19").Value = -1104.5238

# Row 33
$ws.Range("H33").Value = 649.6
$ws.Range("I33").Value = 699.4666999999999
$ws.Range("K33").Value = 699.4666999999999
$ws.Range("M33").Value = -470.4666999999999

# Row 103
$ws.Range("H103").Value = 1653.1177
$ws.Range("I103").Value = 385.42856
$ws.Range("J103").Value = 2540.5
$ws.Range("K103").Value = 1156.28568
$ws.Range("L103").Value = 7621.5
$ws.Range("M103").Value = -570.28568
$ws.Range("N103").Value = -8793.5

# Row 116
$ws.Range("H116").Value = 51206.863
$ws.Range("I116").Value = 73466.92999999999
$ws.Range("J116").Value = 3506.7144
$ws.Range("K116").Value = 73466.92999999999
$ws.Range("L116").Value = 3506.7144
$ws.Range("M116").Value = -70024.92999999999
$ws.Range("N116").Value = -10390.7144

# Row 137
$ws.Range("H137").Value = 5697.788
$ws.Range("I137").Value = 6765.88
$ws.Range("J137").Value = 2360
$ws.Range("K137").Value = 20297.64
$ws.Range("L137").Value = 7080
$ws.Range("M137").Value = -17747.64
$ws.Range("N137").Value = -12180

$ws = $wb.Worksheets.Item("ARM")
# Row 13
$ws.Range("H13").Value = 2750
$ws.Range("J13").Value = 2750
$ws.Range("L13").Value = 2750
$ws.Range("N13").Value = -3038

# Row 61
$ws.Range("H61").Value = 648694.9399999999
$ws.Range("I61").Value = 529212.7
$ws.Range("J61").Value = 837875.2
$ws.Range("K61").Value = 529212.7
$ws.Range("L61").Value = 837875.2
$ws.Range("M61").Value = -529000.7
$ws.Range("N61").Value = -838299.2

# Row 63
$ws.Range("H63").Value = 3475
$ws.Range("I63").Value = 3860
$ws.Range("J63").Value = 2833.3333
$ws.Range("K63").Value = 3860
$ws.Range("L63").Value = 2833.3333
$ws.Range("M63").Value = -3174
$ws.Range("N63").Value = -4205.3333

# Row 66
$ws.Range("H66").Value = 3475
$ws.Range("I66").Value = 3860
$ws.Range("J66").Value = 2833.3333
$ws.Range("K66").Value = 19300
$ws.Range("L66").Value = 14166.6665
$ws.Range("M66").Value = -15868
$ws.Range("N66").Value = -21030.6665

# Row 74
$ws.Range("H74").Value = 241041.84
$ws.Range("I74").Value = 295316.7
$ws.Range("J74").Value = 87263.164
$ws.Range("K74").Value = 295316.7
$ws.Range("L74").Value = 87263.164
$ws.Range("M74").Value = -294442.7
$ws.Range("N74").Value = -89011.164

# Row 77
$ws.Range("H77").Value = 241041.84
$ws.Range("I77").Value = 295316.7
$ws.Range("J77").Value = 87263.164
$ws.Range("K77").Value = 1476583.5
$ws.Range("L77").Value = 436315.82
$ws.Range("M77").Value = -1472215.5
$ws.Range("N77").Value = -445051.82

# Row 88
$ws.Range("H88").Value = 1973.4667
$ws.Range("I88").Value = 1794.1212
$ws.Range("K88").Value = 1794.1212
$ws.Range("M88").Value = -1388.1212

# Row 91
$ws.Range("H91").Value = 1973.4667
$ws.Range("I91").Value = 1794.1212
$ws.Range("K91").Value = 1794.1212
$ws.Range("M91").Value = -390.1212

# Row 132
$ws.Range("H132").Value = 29694.553
$ws.Range("I132").Value = 38732.395
$ws.Range("J132").Value = 4388.6
$ws.Range("K132").Value = 116197.185
$ws.Range("L132").Value = 13165.8
$ws.Range("M132").Value = -113667.185
$ws.Range("N132").Value = -18225.8

# Row 136
$ws.Range("H136").Value = 648694.9399999999
$ws.Range("I136").Value = 529212.7
$ws.Range("J136").Value = 837875.2
$ws.Range("K136").Value = 1587638.1
$ws.Range("L136").Value = 2513625.6
$ws.Range("M136").Value = -1585088.1
$ws.Range("N136").Value = -2518725.6

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 330.2857
$ws.Range("I22").Value = 292
$ws.Range("K22").Value = 292
$ws.Range("M22").Value = -119

# Row 86
$ws.Range("H86").Value = 3907.44
$ws.Range("I86").Value = 7663.25
$ws.Range("J86").Value = 2140
$ws.Range("K86").Value = 7663.25
$ws.Range("L86").Value = 2140
$ws.Range("M86").Value = -6540.25
$ws.Range("N86").Value = -4386

# Row 89
$ws.Range("H89").Value = 3907.44
$ws.Range("I89").Value = 7663.25
$ws.Range("J89").Value = 2140
$ws.Range("K89").Value = 38316.25
$ws.Range("L89").Value = 10700
$ws.Range("M89").Value = -32700.25
$ws.Range("N89").Value = -21932

# Row 99
$ws.Range("H99").Value = 6781.222
$ws.Range("I99").Value = 10824
$ws.Range("J99").Value = 1727.75
$ws.Range("K99").Value = 10824
$ws.Range("L99").Value = 1727.75
$ws.Range("M99").Value = -9326
$ws.Range("N99").Value = -4723.75

# Row 134
$ws.Range("H134").Value = 4879.706
$ws.Range("I134").Value = 5077.3076
$ws.Range("J134").Value = 4237.5
$ws.Range("K134").Value = 15231.9228
$ws.Range("L134").Value = 12712.5
$ws.Range("M134").Value = -12696.9228
$ws.Range("N134").Value = -17782.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2129.5918
$ws.Range("I31").Value = 1224.25
$ws.Range("J31").Value = 6153.3335
$ws.Range("K31").Value = 1224.25
$ws.Range("L31").Value = 6153.3335
$ws.Range("M31").Value = -929.25
$ws.Range("N31").Value = -6743.3335

# Row 34
$ws.Range("H34").Value = 2129.5918
$ws.Range("I34").Value = 1224.25
$ws.Range("J34").Value = 6153.3335
$ws.Range("K34").Value = 1224.25
$ws.Range("L34").Value = 6153.3335
$ws.Range("M34").Value = -1022.25
$ws.Range("N34").Value = -6557.3335

# Row 58
$ws.Range("H58").Value = 7267
$ws.Range("I58").Value = 10279.363
$ws.Range("J58").Value = 3125
$ws.Range("K58").Value = 10279.363
$ws.Range("L58").Value = 3125
$ws.Range("M58").Value = -10076.363
$ws.Range("N58").Value = -3531

# Row 132
$ws.Range("H132").Value = 2315.8857
$ws.Range("I132").Value = 1077
$ws.Range("J132").Value = 5018.909
$ws.Range("K132").Value = 3231
$ws.Range("L132").Value = 15056.727
$ws.Range("M132").Value = -701
$ws.Range("N132").Value = -20116.727

# Row 134
$ws.Range("H134").Value = 2204.889
$ws.Range("I134").Value = 1144.1818
$ws.Range("J134").Value = 3871.7144
$ws.Range("K134").Value = 3432.5454
$ws.Range("L134").Value = 11615.1432
$ws.Range("M134").Value = -897.5454
$ws.Range("N134").Value = -16685.1432

# Row 136
$ws.Range("H136").Value = 7267
$ws.Range("I136").Value = 10279.363
$ws.Range("J136").Value = 3125
$ws.Range("K136").Value = 30838.089
$ws.Range("L136").Value = 9375
$ws.Range("M136").Value = -28288.089
$ws.Range("N136").Value = -14475

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 423.15384
$ws.Range("I11").Value = 300.33334
$ws.Range("J11").Value = 699.5
$ws.Range("K11").Value = 901.0000200000001
$ws.Range("L11").Value = 2098.5
$ws.Range("M11").Value = -761.0000200000001
$ws.Range("N11").Value = -2378.5

# Row 29
$ws.Range("H29").Value = 112.75
$ws.Range("I29").Value = 30.5
$ws.Range("J29").Value = 195
$ws.Range("K29").Value = 91.5
$ws.Range("L29").Value = 585
$ws.Range("M29").Value = 185.5
$ws.Range("N29").Value = -1139

# Row 92
$ws.Range("H92").Value = 22727476
$ws.Range("I92").Value = 27777958
$ws.Range("K92").Value = 83333874
$ws.Range("M92").Value = -83332626

# Row 97
$ws.Range("H97").Value = 701.3333
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 902
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 2706
$ws.Range("M97").Value = -404
$ws.Range("N97").Value = -3698

# Row 131
$ws.Range("H131").Value = 2080.1462
$ws.Range("I131").Value = 3079.0908
$ws.Range("J131").Value = 1713.8667
$ws.Range("K131").Value = 9237.2724
$ws.Range("L131").Value = 5141.6001
$ws.Range("M131").Value = -4197.2724
$ws.Range("N131").Value = -15221.6001

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 5554.9414
$ws.Range("I102").Value = 4083.6365
$ws.Range("J102").Value = 8252.333000000001
$ws.Range("K102").Value = 4083.6365
$ws.Range("L102").Value = 8252.333000000001
$ws.Range("M102").Value = -2461.6365
$ws.Range("N102").Value = -11496.333

# Row 122
$ws.Range("H122").Value = 995.6923
$ws.Range("I122").Value = 759.4
$ws.Range("J122").Value = 1783.3334
$ws.Range("K122").Value = 2278.2
$ws.Range("L122").Value = 5350.0002
$ws.Range("M122").Value = 171.8000000000002
$ws.Range("N122").Value = -10250.0002

# Row 132
$ws.Range("H132").Value = 3869.1428
$ws.Range("I132").Value = 3460.074
$ws.Range("J132").Value = 5249.75
$ws.Range("K132").Value = 10380.222
$ws.Range("L132").Value = 15749.25
$ws.Range("M132").Value = -7850.222
$ws.Range("N132").Value = -20809.25

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 11911801
$ws.Range("I132").Value = 17243578
$ws.Range("J132").Value = 17838.309
$ws.Range("K132").Value = 51730734
$ws.Range("L132").Value = 53514.927
$ws.Range("M132").Value = -51728204
$ws.Range("N132").Value = -58574.927

# Row 136
$ws.Range("H136").Value = 4877.4
$ws.Range("I136").Value = 2768.9092
$ws.Range("J136").Value = 7454.4443
$ws.Range("K136").Value = 8306.7276
$ws.Range("L136").Value = 22363.3329
$ws.Range("M136").Value = -5756.7276
$ws.Range("N136").Value = -27463.3329

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 898.375
$ws.Range("I7").Value = 883.8570999999999
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 883.8570999999999
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -770.8570999999999
$ws.Range("N7").Value = -1226

# Row 122
$ws.Range("H122").Value = 57549.723
$ws.Range("I122").Value = 68446
$ws.Range("J122").Value = 3068.3333
$ws.Range("K122").Value = 205338
$ws.Range("L122").Value = 9204.999899999999
$ws.Range("M122").Value = -202888
$ws.Range("N122").Value = -14104.9999

# Row 132
$ws.Range("H132").Value = 22729072
$ws.Range("I132").Value = 29412828
$ws.Range("J132").Value = 4304.8
$ws.Range("K132").Value = 88238484
$ws.Range("L132").Value = 12914.4
$ws.Range("M132").Value = -88235954
$ws.Range("N132").Value = -17974.4

# Row 136
$ws.Range("H136").Value = 18057738
$ws.Range("I136").Value = 31286232
$ws.Range("J136").Value = 419747.9
$ws.Range("K136").Value = 93858696
$ws.Range("L136").Value = 1259243.7
$ws.Range("M136").Value = -93856146
$ws.Range("N136").Value = -1264343.7
